$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 666.6875
$ws.Range("I33").Value = 732.6429000000001
$ws.Range("K33").Value = 732.6429000000001
$ws.Range("M33").Value = -503.6429000000001
$ws.Range("H62").Value = 6037.2666
$ws.Range("I62").Value = 5548.25
$ws.Range("K62").Value = 5548.25
$ws.Range("M62").Value = -4924.25
$ws.Range("H65").Value = 6037.2666
$ws.Range("I65").Value = 5548.25
$ws.Range("K65").Value = 27741.25
$ws.Range("M65").Value = -24621.25
$ws.Range("H86").Value = 6301
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 6301
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 6301
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -8547
$ws.Range("H89").Value = 6301
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 6301
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 31505
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -42737
$ws.Range("H105").Value = 60000
$ws.Range("J105").Value = 60000
$ws.Range("L105").Value = 60000
$ws.Range("N105").Value = -66988
$ws.Range("H132").Value = 1379.2778
$ws.Range("I132").Value = 1252.625
$ws.Range("K132").Value = 3757.875
$ws.Range("M132").Value = -1227.875
$ws.Range("H138").Value = 4155.68
$ws.Range("J138").Value = 4324.3164
$ws.Range("L138").Value = 12972.9492
$ws.Range("N138").Value = -23252.9492
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1221.5
$ws.Range("I4").Value = 1096.8
$ws.Range("J4").Value = 1429.3334
$ws.Range("K4").Value = 1096.8
$ws.Range("L4").Value = 1429.3334
$ws.Range("M4").Value = -980.8
$ws.Range("N4").Value = -1661.3334
$ws.Range("H5").Value = 303.83334
$ws.Range("I5").Value = 318.6
$ws.Range("J5").Value = 230
$ws.Range("K5").Value = 318.6
$ws.Range("L5").Value = 230
$ws.Range("M5").Value = -206.6
$ws.Range("N5").Value = -454
$ws.Range("H6").Value = 8049.5
$ws.Range("I6").Value = 1099
$ws.Range("K6").Value = 1099
$ws.Range("M6").Value = -926
$ws.Range("H37").Value = 12074.5
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 40000
$ws.Range("N37").Value = -40546
$ws.Range("H61").Value = 7086.294
$ws.Range("I61").Value = 5784.2856
$ws.Range("K61").Value = 5784.2856
$ws.Range("M61").Value = -5572.2856
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H96").Value = 18498.428
$ws.Range("J96").Value = 18498.428
$ws.Range("L96").Value = 18498.428
$ws.Range("N96").Value = -23990.428
$ws.Range("H101").Value = 46600
$ws.Range("J101").Value = 46600
$ws.Range("L101").Value = 46600
$ws.Range("N101").Value = -53090
$ws.Range("H136").Value = 7086.294
$ws.Range("I136").Value = 5784.2856
$ws.Range("K136").Value = 17352.8568
$ws.Range("M136").Value = -14802.8568
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 303.83334
$ws.Range("I4").Value = 318.6
$ws.Range("J4").Value = 230
$ws.Range("K4").Value = 318.6
$ws.Range("L4").Value = 230
$ws.Range("M4").Value = -203.6
$ws.Range("N4").Value = -460
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("K19").Value = 1000
$ws.Range("M19").Value = -827
$ws.Range("H22").Value = 1099.6666
$ws.Range("I22").Value = 1114.4286
$ws.Range("J22").Value = 1079
$ws.Range("K22").Value = 1114.4286
$ws.Range("L22").Value = 1079
$ws.Range("M22").Value = -941.4286
$ws.Range("N22").Value = -1425
$ws.Range("H105").Value = 2649
$ws.Range("I105").Value = 1473.5
$ws.Range("K105").Value = 1473.5
$ws.Range("M105").Value = 273.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19840.637
$ws.Range("I31").Value = 32651.5
$ws.Range("K31").Value = 32651.5
$ws.Range("M31").Value = -32356.5
$ws.Range("H34").Value = 19840.637
$ws.Range("I34").Value = 32651.5
$ws.Range("K34").Value = 32651.5
$ws.Range("M34").Value = -32449.5
$ws.Range("H58").Value = 1400
$ws.Range("I58").Value = 1400
$ws.Range("K58").Value = 1400
$ws.Range("M58").Value = -1197
$ws.Range("H86").Value = 8976.75
$ws.Range("I86").Value = 11474.5
$ws.Range("J86").Value = 6479
$ws.Range("K86").Value = 11474.5
$ws.Range("L86").Value = 6479
$ws.Range("M86").Value = -10351.5
$ws.Range("N86").Value = -8725
$ws.Range("H89").Value = 8976.75
$ws.Range("I89").Value = 11474.5
$ws.Range("J89").Value = 6479
$ws.Range("K89").Value = 57372.5
$ws.Range("L89").Value = 32395
$ws.Range("M89").Value = -51756.5
$ws.Range("N89").Value = -43627
$ws.Range("H107").Value = 849.4286
$ws.Range("I107").Value = 379.1111
$ws.Range("J107").Value = 1696
$ws.Range("K107").Value = 379.1111
$ws.Range("L107").Value = 1696
$ws.Range("M107").Value = 1540.8889
$ws.Range("N107").Value = -5536
$ws.Range("H132").Value = 2643.5
$ws.Range("J132").Value = 2763.3333
$ws.Range("L132").Value = 8289.999899999999
$ws.Range("N132").Value = -13349.9999
$ws.Range("H134").Value = 3512.8572
$ws.Range("I134").Value = 3558
$ws.Range("J134").Value = 3400
$ws.Range("K134").Value = 10674
$ws.Range("L134").Value = 10200
$ws.Range("M134").Value = -8139
$ws.Range("N134").Value = -15270
$ws.Range("H136").Value = 1400
$ws.Range("I136").Value = 1400
$ws.Range("K136").Value = 4200
$ws.Range("M136").Value = -1650
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 10583.9
$ws.Range("I70").Value = 2010
$ws.Range("J70").Value = 16299.833
$ws.Range("K70").Value = 6030
$ws.Range("L70").Value = 48899.499
$ws.Range("M70").Value = -5715
$ws.Range("N70").Value = -49529.499
$ws.Range("H73").Value = 10583.9
$ws.Range("I73").Value = 2010
$ws.Range("J73").Value = 16299.833
$ws.Range("K73").Value = 6030
$ws.Range("L73").Value = 48899.499
$ws.Range("M73").Value = -4938
$ws.Range("N73").Value = -51083.499
$ws.Range("H107").Value = 239.5
$ws.Range("I107").Value = 165.4
$ws.Range("J107").Value = 363
$ws.Range("K107").Value = 496.2
$ws.Range("L107").Value = 1089
$ws.Range("M107").Value = 1423.8
$ws.Range("N107").Value = -4929
$ws.Range("H132").Value = 6833.3335
$ws.Range("I132").Value = 8250
$ws.Range("K132").Value = 74250
$ws.Range("M132").Value = -71720
$ws.Range("H134").Value = 394.5
$ws.Range("I134").Value = 394.5
$ws.Range("K134").Value = 1183.5
$ws.Range("M134").Value = 3886.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 30247.5
$ws.Range("J39").Value = 30247.5
$ws.Range("L39").Value = 30247.5
$ws.Range("N39").Value = -31311.5
$ws.Range("H80").Value = 13186.091
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 14204.7
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 14204.7
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -16200.7
$ws.Range("H83").Value = 13186.091
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 14204.7
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 71023.5
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -81007.5
$ws.Range("H121").Value = 48896.5
$ws.Range("J121").Value = 48896.5
$ws.Range("L121").Value = 48896.5
$ws.Range("N121").Value = -52390.5
$ws.Range("H122").Value = 2401.5
$ws.Range("I122").Value = 1803
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5409
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2959
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 1689.5
$ws.Range("I132").Value = 1528.9
$ws.Range("K132").Value = 4586.700000000001
$ws.Range("M132").Value = -2056.700000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 40514.832
$ws.Range("I74").Value = 36297.75
$ws.Range("K74").Value = 36297.75
$ws.Range("M74").Value = -35299.75
$ws.Range("H77").Value = 40514.832
$ws.Range("I77").Value = 36297.75
$ws.Range("K77").Value = 108893.25
$ws.Range("M77").Value = -103901.25
$ws.Range("H95").Value = 44326.668
$ws.Range("J95").Value = 44326.668
$ws.Range("L95").Value = 44326.668
$ws.Range("N95").Value = -49818.668
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 53666.668
$ws.Range("J92").Value = 53666.668
$ws.Range("L92").Value = 53666.668
$ws.Range("N92").Value = -58658.668
$ws.Range("H95").Value = 16666.5
$ws.Range("J95").Value = 15374.75
$ws.Range("L95").Value = 15374.75
$ws.Range("N95").Value = -20866.75
$ws.Range("H132").Value = 7146
$ws.Range("I132").Value = 7146
$ws.Range("K132").Value = 21438
$ws.Range("M132").Value = -18908
$ws.Range("H136").Value = 4635.75
$ws.Range("I136").Value = 5947.6665
$ws.Range("J136").Value = 700
$ws.Range("K136").Value = 17842.9995
$ws.Range("L136").Value = 2100
$ws.Range("M136").Value = -15292.9995
$ws.Range("N136").Value = -7200
